$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("43:45").Insert()
Write-Host "inserted"
$v = $ws.Range("B44").Value()
Write-Host ("B44 after insert: " + $v)
$v2 = $ws.Range("B47").Value()
Write-Host ("B47 after insert: " + $v2)
